$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NewImportLogic_1 - Test_Automation_1"
$ws.Range("C2").Value = "NewImportLogic_1 - Test_Automation_1_radio_button"

$ws.Range("D2").Style = "Normal"

$ws.Range("D2").Select()
